$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 data entry (Data / Dinheiro / Cartao / Pix) ---
$ws.Range("A2").Value = 45770.806770833333
$ws.Range("B2").Value = 150
$ws.Range("C2").Value = 1492.5
$ws.Range("D2").Value = 399.5

# --- Column A sized to fit the new date/time value ---
$ws.Columns.Item(1).ColumnWidth = 16.8

# --- Selection moves to C3 ---
$ws.Range("C3").Select() | Out-Null
